$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(84, 4, 13, 101, 0.9545454545454546, 0.865979381443299)
    3 = @(68, 10, 30, 108, 0.8717948717948718, 0.6938775510204082)
    4 = @(41, 16, 55, 112, 0.7192982456140351, 0.4270833333333333)
    5 = @(73, 7, 25, 105, 0.9125, 0.7448979591836735)
    6 = @(84, 6, 13, 103, 0.9333333333333333, 0.865979381443299)
    7 = @(87, 2, 12, 101, 0.9775280898876404, 0.8787878787878788)
    8 = @(83, 5, 15, 103, 0.9431818181818182, 0.8469387755102041)
    9 = @(57, 25, 36, 118, 0.6951219512195121, 0.6129032258064516)
    10 = @(38, 11, 57, 106, 0.7755102040816326, 0.4)
    11 = @(85, 6, 14, 105, 0.9340659340659341, 0.8585858585858586)
    12 = @(79, 9, 19, 107, 0.8977272727272727, 0.8061224489795918)
    13 = @(93, 5, 6, 104, 0.9489795918367347, 0.9393939393939394)
    14 = @(66, 4, 32, 102, 0.9428571428571428, 0.673469387755102)
    15 = @(75, 5, 24, 104, 0.9375, 0.7575757575757576)
    16 = @(21, 1, 30, 52, 0.9545454545454546, 0.4117647058823529)
    17 = @(83, 4, 15, 102, 0.9540229885057471, 0.8469387755102041)
    18 = @(96, 0, 2, 98, 1, 0.9795918367346939)
    19 = @(72, 5, 23, 100, 0.935064935064935, 0.7578947368421053)
    20 = @(66, 11, 32, 109, 0.8571428571428571, 0.673469387755102)
    21 = @(59, 13, 38, 110, 0.8194444444444444, 0.6082474226804123)
    22 = @(44, 8, 50, 102, 0.8461538461538461, 0.4680851063829787)
    23 = @(86, 1, 12, 99, 0.9885057471264368, 0.8775510204081632)
    24 = @(87, 4, 12, 103, 0.9560439560439561, 0.8787878787878788)
    25 = @(83, 7, 16, 106, 0.9222222222222223, 0.8383838383838383)
    26 = @(58, 8, 38, 104, 0.8787878787878788, 0.6041666666666666)
    27 = @(69, 19, 22, 110, 0.7840909090909091, 0.7582417582417582)
    28 = @(74, 13, 15, 102, 0.8505747126436781, 0.8314606741573034)
    29 = @(44, 23, 37, 104, 0.6567164179104478, 0.5432098765432098)
    30 = @(84, 4, 15, 103, 0.9545454545454546, 0.8484848484848485)
    31 = @(78, 6, 18, 102, 0.9285714285714286, 0.8125)
    32 = @(61, 12, 31, 104, 0.8356164383561644, 0.6630434782608695)
    33 = @(93, 0, 6, 99, 1, 0.9393939393939394)
    34 = @(73, 3, 24, 100, 0.9605263157894737, 0.7525773195876289)
    35 = @(90, 1, 9, 100, 0.9890109890109891, 0.9090909090909091)
    36 = @(40, 19, 58, 117, 0.6779661016949152, 0.4081632653061225)
    37 = @(60, 12, 33, 105, 0.8333333333333334, 0.6451612903225806)
    38 = @(77, 7, 21, 105, 0.9166666666666666, 0.7857142857142857)
    39 = @(28, 23, 61, 112, 0.5490196078431373, 0.3146067415730337)
    40 = @(80, 4, 19, 103, 0.9523809523809523, 0.8080808080808081)
    41 = @(36, 7, 12, 55, 0.8372093023255814, 0.75)
    42 = @(88, 0, 9, 97, 1, 0.9072164948453608)
    43 = @(64, 8, 28, 100, 0.8888888888888888, 0.6956521739130435)
    44 = @(91, 2, 8, 101, 0.978494623655914, 0.9191919191919192)
    45 = @(89, 6, 9, 104, 0.9368421052631579, 0.9081632653061225)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item([int]$row, 3).Value = $vals[0]  # C: TP
    $ws.Cells.Item([int]$row, 4).Value = $vals[1]  # D: FP
    $ws.Cells.Item([int]$row, 5).Value = $vals[2]  # E: FN
    $ws.Cells.Item([int]$row, 6).Value = $vals[3]  # F: Total
    $ws.Cells.Item([int]$row, 7).Value = $vals[4]  # G: Precision
    $ws.Cells.Item([int]$row, 8).Value = $vals[5]  # H: Recall
}
